# Apply the "added to LFO section" edit.
#
# Net content changes (per the commit's diff):
#   1. "The MIDI to Gate/Trigger is a simple device..."   -> "The MindBrain MIDI to Gate/Trigger is a simple device..."
#   2. "Text Here "  (LFO placeholder body paragraph)      -> full LFO description paragraph
#   3. "The Dual Envelope Generator is a versatile..."     -> "The MindBrain Dual Envelope Generator is a versatile..."
#   4. "The Dual Random Generator device can be used..."   -> "The MindBrain Dual Random Generator device can be used..."

$d = $word.ActiveDocument

# 1) MIDI to Gate/Trigger intro paragraph: insert "MindBrain " after "The ".
$d.Content.Find.Execute(
    "The MIDI to Gate/Trigger is a simple device",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The MindBrain MIDI to Gate/Trigger is a simple device",
    2) | Out-Null

# 2) Replace the "Text Here " placeholder under the LFO heading with the real
#    device description.
$d.Content.Find.Execute(
    "Text Here ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The MindBrain Dual Low Frequency Oscillator device is a simple but deep LFO module. Each of its two sides can be set to one of seven different wave shapes – Sine, Ramp Up, Ramp Down, Triangle, Rectangle, Random and Bin. The LFO can be set to sync with Live’s clock, or can be set to an independent frequency up to 2KHz.",
    2) | Out-Null

# 3) Dual Envelope Generator intro paragraph: insert "MindBrain " after "The ".
$d.Content.Find.Execute(
    "The Dual Envelope Generator is a versatile",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The MindBrain Dual Envelope Generator is a versatile",
    2) | Out-Null

# 4) Dual Random Generator intro paragraph: insert "MindBrain " after "The ".
$d.Content.Find.Execute(
    "The Dual Random Generator device can be used",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The MindBrain Dual Random Generator device can be used",
    2) | Out-Null

Write-Output "Edits applied"
